# Apply cryptocurrency price/volume updates to match the source diff.
# Column D (Price) values that look numeric are prefixed with a leading
# apostrophe so Excel stores them as literal text (preserving exact
# formatting such as trailing zeros and thousand-dot groupings) instead
# of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.099.46"
$ws.Range("E2").Value = "  -2.68%  "

$ws.Range("D3").Value = "1.844.15"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'0.6905"
$ws.Range("E5").Value = "  -6.74%  "

$ws.Range("D6").Value = "'236.97"
$ws.Range("E6").Value = "  -2.25%  "

$ws.Range("D8").Value = "'0.3036"
$ws.Range("E8").Value = "  -3.65%  "

$ws.Range("D9").Value = "'0.07503"
$ws.Range("E9").Value = "  +3.92%  "

$ws.Range("D10").Value = "'23.31"
$ws.Range("E10").Value = "  -5.50%  "

$ws.Range("D11").Value = "'0.08083"
$ws.Range("E11").Value = "  -2.86%  "

$ws.Range("D12").Value = "1.842.62"
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").Value = "'0.7204"
$ws.Range("E13").Value = "  -3.91%  "

$ws.Range("D14").Value = "'5.164"
$ws.Range("E14").Value = "  -4.01%  "

$ws.Range("D15").Value = "'88.63"
$ws.Range("E15").Value = "  -3.84%  "

$ws.Range("D16").Value = "29.191.28"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "'5.773"
$ws.Range("E17").Value = "  -5.45%  "

$ws.Range("D18").Value = "'241.21"
$ws.Range("E18").Value = "  -2.31%  "

$ws.Range("D19").Value = "'0.000007657"
$ws.Range("E19").Value = "  -2.17%  "

$ws.Range("D20").Value = "'12.97"
$ws.Range("E20").Value = "  -4.21%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "2.106.10"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'7.606"
$ws.Range("E24").Value = "  -4.81%  "

$ws.Range("D25").Value = "'8.983"
$ws.Range("E25").Value = "  -3.25%  "

$ws.Range("E26").Value = "  -5.32%  "

$ws.Range("D27").Value = "'161.32"

$ws.Range("D28").Value = "'18.02"
$ws.Range("E28").Value = "  -3.33%  "

$ws.Range("D29").Value = "'1.925"
$ws.Range("E29").Value = "  -4.67%  "

$ws.Range("D30").Value = "'1.378"
$ws.Range("E30").Value = "  -7.76%  "

$ws.Range("D31").Value = "'4.428"
$ws.Range("E31").Value = "  -3.19%  "

$ws.Range("D32").Value = "'1.487"
$ws.Range("E32").Value = "  -3.07%  "

$ws.Range("D33").Value = "'4.030"
$ws.Range("E33").Value = "  -4.43%  "

$ws.Range("D34").Value = "'0.05190"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  -4.79%  "

$ws.Range("D36").Value = "'0.7107"
$ws.Range("E36").Value = "  -5.12%  "

$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Value = "'2.654"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").Value = "'0.01856"
$ws.Range("E39").Value = "  -5.32%  "

$ws.Range("D40").Value = "'2.673"
$ws.Range("E40").Value = "  -2.88%  "

$ws.Range("D41").Value = "'0.9159"
$ws.Range("E41").Value = "  +6.06%  "

$ws.Range("D42").Value = "'5.924"
$ws.Range("E42").Value = "  -3.30%  "

$ws.Range("D43").Value = "'0.4266"
$ws.Range("E43").Value = "  -5.88%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.049.92"
$ws.Range("E44").Value = "  -6.10%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'69.71"
$ws.Range("E45").Value = "  -3.60%  "

$ws.Range("D46").Value = "'1.000"

$ws.Range("D47").Value = "'102.36"
$ws.Range("E47").Value = "  -1.95%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.014.57"
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.146"
$ws.Range("E49").Value = "  -6.19%  "

$ws.Range("D50").Value = "'1.740"
$ws.Range("E50").Value = "  -6.44%  "

$ws.Range("D51").Value = "'9.217"
$ws.Range("E51").Value = "  -2.88%  "
